# "End of section 3" -- add the Migrations/DataBase rows (3-6) to the
# "Visual code commands" sheet, add the "additional comments" column (D),
# switch the active tab to "Visual code commands", and move the selection
# left on "CMD commands".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CMD commands")
$ws2 = $wb.Worksheets.Item("Visual code commands")

# --- Fill in the new cells in the exact order the original author typed
# them (this reproduces the shared-string table ordering in the saved
# file: row 6 first, then the new D1 header, then rows 3-5). ---

$ws2.Range("C6").Value = "dotnet ef database drop -p Infrastructure -s API"
$ws2.Range("B6").Value = "Drop Database when we have 3 separate projects"

$ws2.Range("D1").Value = "additional comments"

$ws2.Range("D3").Value = "p = project  (where Dbcontext in)  s = start (where startup file in)"
$ws2.Range("C3").Value = "dotnet ef migrations remove -p Infrastructure -s API"

$ws2.Range("A5").Value = "DataBase"

$ws2.Range("B3").Value = "Remove migrations"
$ws2.Range("B4").Value = "Add migrations"
$ws2.Range("C4").Value = "dotnet ef migrations add InitialCreate -p Infrastructure -s API -o Data/Migrations"
$ws2.Range("D4").Value = "p = project  (where Dbcontext in)  s = start (where startup file in) o = path to create the migration"
$ws2.Range("C5").Value = "dotnet ef database update -p Infrastructure -s API"
$ws2.Range("B5").Value = "Update database"

# --- Fill in the remaining "anchor" column A cells + D6 ---
$ws2.Range("A3").Value = "Migrations"
$ws2.Range("A4").Value = "Migrations"
$ws2.Range("A6").Value = "DataBase"
$ws2.Range("D6").Value = "p = project  (where Dbcontext in)  s = start (where startup file in)"

# --- Column widths (C widened / loses best-fit, new D column best-fit) ---
$ws2.Columns.Item(3).ColumnWidth = 82.83333333333333
$ws2.Columns.Item(4).ColumnWidth = 87.5

# --- Page setup: portrait orientation on sheet 2 ---
$ws2.PageSetup.Orientation = 1

# --- Selection / active-tab bookkeeping ---
# "CMD commands" selection moves to A14 and stops being the active tab.
$ws1.Range("A14").Select() | Out-Null

# "Visual code commands" becomes the active tab with C4 selected.
$ws2.Activate() | Out-Null
$ws2.Range("C4").Select() | Out-Null
